# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the rest of the header row (bold, centered, bordered)
# by copying the formatting from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record for every player row (rows 2 through 48).
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 97   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 65   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
